$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '97.290.29'
$ws.Range('E2').Value = '  +2.26%  '

# Row 3
$ws.Range('D3').Value = '3.578.77'
$ws.Range('E3').Value = '  +0.14%  '

# Row 4
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.01'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.46%  '

# Row 6
$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.73'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +16.95%  '

# Row 7
$ws.Range('B7').Value = 'BNB'
$ws.Range('C7').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '655.21'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.37%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.430'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +7.84%  '

# Row 9
$ws.Range('E9').Value = '  -0.12%  '

# Row 10
$ws.Range('E10').Value = '  +4.78%  '

# Row 11
$ws.Range('D11').Value = '3.575.08'
$ws.Range('E11').Value = '  +0.09%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '44.26'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +4.52%  '

# Row 13
$ws.Range('E13').Value = '  +0.46%  '

# Row 14
$ws.Range('E14').Value = '  -0.50%  '

# Row 15
$ws.Range('D15').Value = '4.242.84'
$ws.Range('E15').Value = '  +0.04%  '

# Row 16
$ws.Range('D16').Value = '96.974.29'
$ws.Range('E16').Value = '  +2.07%  '

# Row 17
$ws.Range('E17').Value = '  +3.32%  '

# Row 18
$ws.Range('E18').Value = '  +12.01%  '

# Row 19
$ws.Range('D19').Value = '3.578.37'
$ws.Range('E19').Value = '  +0.07%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.68'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.58%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.03'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.24%  '

# Row 23
$ws.Range('E23').Value = '  +1.06%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '514.64'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.21%  '

# Row 25
$ws.Range('E25').Value = '  +5.59%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.91'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.95%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '101.83'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +7.02%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '13.09'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.42%  '

# Row 29
$ws.Range('D29').Value = '3.769.97'
$ws.Range('E29').Value = '  +0.12%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.168'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +16.96%  '

# Row 31
$ws.Range('E31').Value = '  -1.49%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.99'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.17%  '

# Row 33
$ws.Range('E33').Value = '  -0.04%  '

# Row 34
$ws.Range('E34').Value = '  +4.62%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.05%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.90'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.01%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.79'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +4.13%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '613.20'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +5.76%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.567'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.70%  '

# Row 40
$ws.Range('E40').Value = '  -1.29%  '

# Row 41
$ws.Range('E41').Value = '  +2.65%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.94'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +6.59%  '

# Row 43
$ws.Range('E43').Value = '  -0.05%  '

# Row 44
$ws.Range('E44').Value = '  +2.10%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.02'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +5.04%  '

# Row 46
$ws.Range('E46').Value = '  +6.25%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.33'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.33%  '

# Row 48
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.418'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +35.91%  '

# Row 49
$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.61'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.94%  '

# Row 50
$ws.Range('E50').Value = '  +4.49%  '

# Row 51
$ws.Range('E51').Value = '  +7.72%  '
